$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 309). The value was bumped from 46074 to 46075
# (one day later) for all of them.
$ws.Range("C2:C309").Value = 46075
